$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "followers" column (L) for rows 2-17: scraped value refreshed to 14,000,000
$ws.Range("L2:L17").Value = 14000000

# Column widths picked up while reviewing/adjusting the newly scraped data
$ws.Columns("B").ColumnWidth = 10.85546875
$ws.Columns("L").ColumnWidth = 11

# Leave the selection on the last edited cell, matching where the author ended up
$ws.Range("L17").Select() | Out-Null
